$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1275.4231
$ws.Range("J17").Value = 1225.7858
$ws.Range("L17").Value = 3677.3574
$ws.Range("N17").Value = -4013.3574
$ws.Range("H33").Value = 124.76
$ws.Range("I33").Value = 71.38095
$ws.Range("K33").Value = 71.38095
$ws.Range("M33").Value = 157.61905
$ws.Range("H51").Value = 4879.9
$ws.Range("J51").Value = 5374.875
$ws.Range("L51").Value = 5374.875
$ws.Range("N51").Value = -6342.875
$ws.Range("H62").Value = 1548
$ws.Range("I62").Value = 1500.8
$ws.Range("J62").Value = 1626.6666
$ws.Range("K62").Value = 1500.8
$ws.Range("L62").Value = 1626.6666
$ws.Range("M62").Value = -876.8
$ws.Range("N62").Value = -2874.6666
$ws.Range("H65").Value = 1548
$ws.Range("I65").Value = 1500.8
$ws.Range("J65").Value = 1626.6666
$ws.Range("K65").Value = 7504
$ws.Range("L65").Value = 8133.333000000001
$ws.Range("M65").Value = -4384
$ws.Range("N65").Value = -14373.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8209.943
$ws.Range("I32").Value = 6933.5312
$ws.Range("J32").Value = 10154.952
$ws.Range("K32").Value = 6933.5312
$ws.Range("L32").Value = 10154.952
$ws.Range("M32").Value = -6646.5312
$ws.Range("N32").Value = -10728.952
$ws.Range("H63").Value = 6025278
$ws.Range("I63").Value = 9895964
$ws.Range("J63").Value = 4211.1113
$ws.Range("K63").Value = 9895964
$ws.Range("L63").Value = 4211.1113
$ws.Range("M63").Value = -9895278
$ws.Range("N63").Value = -5583.1113
$ws.Range("H66").Value = 6025278
$ws.Range("I66").Value = 9895964
$ws.Range("J66").Value = 4211.1113
$ws.Range("K66").Value = 49479820
$ws.Range("L66").Value = 21055.5565
$ws.Range("M66").Value = -49476388
$ws.Range("N66").Value = -27919.5565
$ws.Range("H74").Value = 3764.6316
$ws.Range("I74").Value = 4019.6553
$ws.Range("K74").Value = 4019.6553
$ws.Range("M74").Value = -3145.6553
$ws.Range("H77").Value = 3764.6316
$ws.Range("I77").Value = 4019.6553
$ws.Range("K77").Value = 20098.2765
$ws.Range("M77").Value = -15730.2765
$ws.Range("H133").Value = 35890
$ws.Range("J133").Value = 35890
$ws.Range("L133").Value = 35890
$ws.Range("N133").Value = -40950
$ws.Range("H137").Value = 39757.555
$ws.Range("J137").Value = 39757.555
$ws.Range("L137").Value = 39757.555
$ws.Range("N137").Value = -49957.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 811.96
$ws.Range("I94").Value = 768.8421
$ws.Range("J94").Value = 948.5
$ws.Range("K94").Value = 768.8421
$ws.Range("L94").Value = 948.5
$ws.Range("M94").Value = -317.8421
$ws.Range("N94").Value = -1850.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2258.7454
$ws.Range("I58").Value = 1875.3405
$ws.Range("J58").Value = 4511.25
$ws.Range("K58").Value = 1875.3405
$ws.Range("L58").Value = 4511.25
$ws.Range("M58").Value = -1672.3405
$ws.Range("N58").Value = -4917.25
$ws.Range("H107").Value = 656.0625
$ws.Range("I107").Value = 280.1
$ws.Range("K107").Value = 280.1
$ws.Range("M107").Value = 1639.9
$ws.Range("H132").Value = 3980.2
$ws.Range("I132").Value = 1600.8
$ws.Range("J132").Value = 6359.6
$ws.Range("K132").Value = 4802.4
$ws.Range("L132").Value = 19078.8
$ws.Range("M132").Value = -2272.4
$ws.Range("N132").Value = -24138.8
$ws.Range("H136").Value = 2258.7454
$ws.Range("I136").Value = 1875.3405
$ws.Range("J136").Value = 4511.25
$ws.Range("K136").Value = 5626.0215
$ws.Range("L136").Value = 13533.75
$ws.Range("M136").Value = -3076.0215
$ws.Range("N136").Value = -18633.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 217
$ws.Range("J23").Value = 292.1
$ws.Range("L23").Value = 876.3000000000001
$ws.Range("N23").Value = -1346.3
$ws.Range("H61").Value = 268.4
$ws.Range("I61").Value = 85.5
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 256.5
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -41.5
$ws.Range("N61").Value = -3430
$ws.Range("H87").Value = 5789.7144
$ws.Range("I87").Value = 2985.6
$ws.Range("J87").Value = 12800
$ws.Range("K87").Value = 8956.8
$ws.Range("L87").Value = 38400
$ws.Range("M87").Value = -7708.799999999999
$ws.Range("N87").Value = -40896
$ws.Range("H90").Value = 5789.7144
$ws.Range("I90").Value = 2985.6
$ws.Range("J90").Value = 12800
$ws.Range("K90").Value = 26870.4
$ws.Range("L90").Value = 115200
$ws.Range("M90").Value = -20630.4
$ws.Range("N90").Value = -127680
$ws.Range("H103").Value = 1332.5714
$ws.Range("I103").Value = 433.33334
$ws.Range("J103").Value = 2007
$ws.Range("K103").Value = 1300.00002
$ws.Range("L103").Value = 6021
$ws.Range("M103").Value = -421.0000199999999
$ws.Range("N103").Value = -7779
$ws.Range("H113").Value = 653.5952
$ws.Range("I113").Value = 584.2778
$ws.Range("J113").Value = 705.5833
$ws.Range("K113").Value = 1752.8334
$ws.Range("L113").Value = 2116.7499
$ws.Range("M113").Value = 417.1666
$ws.Range("N113").Value = -6456.7499
$ws.Range("H131").Value = 8334244
$ws.Range("I131").Value = 166666930
$ws.Range("J131").Value = 944.6667
$ws.Range("K131").Value = 500000790
$ws.Range("L131").Value = 2834.0001
$ws.Range("M131").Value = -499995750
$ws.Range("N131").Value = -12914.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 8000
$ws.Range("J98").Value = 8000
$ws.Range("L98").Value = 8000
$ws.Range("N98").Value = -13990
$ws.Range("H122").Value = 3115.65
$ws.Range("I122").Value = 1487
$ws.Range("K122").Value = 4461
$ws.Range("M122").Value = -2011
$ws.Range("H132").Value = 3216.7666
$ws.Range("I132").Value = 2357.4546
$ws.Range("J132").Value = 3714.2632
$ws.Range("K132").Value = 7072.3638
$ws.Range("L132").Value = 11142.7896
$ws.Range("M132").Value = -4542.3638
$ws.Range("N132").Value = -16202.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5947.615
$ws.Range("I122").Value = 3534
$ws.Range("K122").Value = 10602
$ws.Range("M122").Value = -8152
$ws.Range("H132").Value = 4261.5425
$ws.Range("I132").Value = 1708.8518
$ws.Range("K132").Value = 5126.555399999999
$ws.Range("M132").Value = -2596.555399999999
$ws.Range("H136").Value = 3982.4211
$ws.Range("I136").Value = 1666.6154
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 4999.8462
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -2449.8462
$ws.Range("N136").Value = -32100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 47637796
$ws.Range("I132").Value = 100000
$ws.Range("J132").Value = 55560760
$ws.Range("K132").Value = 300000
$ws.Range("L132").Value = 166682280
$ws.Range("M132").Value = -297470
$ws.Range("N132").Value = -166687340
